# Update "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - normalize the "Week" labels (W01 -> W1, etc.)
#  - correct two MyForecast values (row 2: 46 -> 47, row 5: 59 -> 58)
#  - store is_holiday_week as a real boolean
# Also refresh the dependent "Min Forecast" value on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before ASIN (old column B) to hold Week_Start_Date.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "Week_Start_Date"

$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)
$myForecast = @(47,60,65,58,58,72,68,68,64,65,62,62,61,55,65,56)

for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2

    $ws.Range("A$r").Value = $weekLabels[$i]

    # Keep the date as literal text (matches source file: inline string, not a date serial).
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $weekStartDates[$i]

    $ws.Range("D$r").Value = $myForecast[$i]

    # is_holiday_week becomes a genuine boolean column.
    $ws.Range("J$r").Value = $false
}

# Summary sheet: "Min Forecast" mirrors the corrected row-2 MyForecast value.
# (Kept as literal text, matching the rest of the Metric/Value column.)
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B14").NumberFormat = "@"
$summary.Range("B14").Value = "47"
